$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values in column B (weekly case counts)
$ws.Range("B4").Value = 436
$ws.Range("B5").Value = 637
$ws.Range("B9").Value = 511
$ws.Range("B12").Value = 471
$ws.Range("B14").Value = 451
$ws.Range("B18").Value = 456
$ws.Range("B21").Value = 304
$ws.Range("B26").Value = 352
$ws.Range("B27").Value = 274
$ws.Range("B29").Value = 305
$ws.Range("B30").Value = 342
$ws.Range("B36").Value = 413
$ws.Range("B37").Value = 449
$ws.Range("B41").Value = 481
$ws.Range("B45").Value = 468

# Add new rows 46 and 47 for weeks 45 and 46
$ws.Range("A46").Value = 45
$ws.Range("B46").Value = 387

$ws.Range("A47").Value = 46
$ws.Range("B47").Value = 66
